{"js": "// Update the worksheet date and each \"a\u00f7b=c, d\" answer cell to the new\n// values described by the commit diff. Every <w:t> run in this document is\n// unique, so a matchCase/matchWholeWord search-and-replace per pair is\n// unambiguous and order-independent.\nconst replacements = [\n  [\"2026-01-27 Tuesday\", \"2026-01-28 Wednesday\"],\n  [\"441\u00f72=220, 1\", \"143\u00f78=17, 7\"],\n  [\"145\u00f72=72, 1\", \"808\u00f73=269, 1\"],\n  [\"396\u00f73=132, 0\", \"342\u00f78=42, 6\"],\n  [\"866\u00f76=144, 2\", \"308\u00f75=61, 3\"],\n  [\"302\u00f73=100, 2\", \"140\u00f76=23, 2\"],\n  [\"859\u00f75=171, 4\", \"990\u00f79=110, 0\"],\n  [\"838\u00f79=93, 1\", \"784\u00f75=156, 4\"],\n  [\"572\u00f76=95, 2\", \"662\u00f79=73, 5\"],\n  [\"652\u00f78=81, 4\", \"947\u00f78=118, 3\"],\n  [\"368\u00f78=46, 0\", \"214\u00f78=26, 6\"],\n  [\"833\u00f76=138, 5\", \"236\u00f78=29, 4\"],\n  [\"793\u00f72=396, 1\", \"508\u00f72=254, 0\"],\n  [\"150\u00f72=75, 0\", \"899\u00f74=224, 3\"],\n  [\"116\u00f77=16, 4\", \"482\u00f76=80, 2\"],\n  [\"916\u00f79=101, 7\", \"437\u00f77=62, 3\"],\n  [\"479\u00f76=79, 5\", \"291\u00f73=97, 0\"],\n  [\"651\u00f72=325, 1\", \"299\u00f74=74, 3\"],\n  [\"530\u00f79=58, 8\", \"667\u00f74=166, 3\"],\n  [\"130\u00f74=32, 2\", \"292\u00f76=48, 4\"],\n  [\"889\u00f78=111, 1\", \"203\u00f72=101, 1\"],\n  [\"331\u00f73=110, 1\", \"882\u00f78=110, 2\"],\n  [\"523\u00f75=104, 3\", \"731\u00f77=104, 3\"],\n  [\"519\u00f73=173, 0\", \"589\u00f77=84, 1\"],\n  [\"550\u00f74=137, 2\", \"201\u00f77=28, 5\"],\n  [\"290\u00f75=58, 0\", \"957\u00f76=159, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and each \"a\u00f7b=c, d\" answer cell to the new\n# values described by the commit diff. Every run of text in this document is\n# unique, so a MatchCase Find/Replace per pair is unambiguous and\n# order-independent.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @('2026-01-27 Tuesday', '2026-01-28 Wednesday'),\n  @('441\u00f72=220, 1', '143\u00f78=17, 7'),\n  @('145\u00f72=72, 1', '808\u00f73=269, 1'),\n  @('396\u00f73=132, 0', '342\u00f78=42, 6'),\n  @('866\u00f76=144, 2', '308\u00f75=61, 3'),\n  @('302\u00f73=100, 2', '140\u00f76=23, 2'),\n  @('859\u00f75=171, 4', '990\u00f79=110, 0'),\n  @('838\u00f79=93, 1', '784\u00f75=156, 4'),\n  @('572\u00f76=95, 2', '662\u00f79=73, 5'),\n  @('652\u00f78=81, 4', '947\u00f78=118, 3'),\n  @('368\u00f78=46, 0', '214\u00f78=26, 6'),\n  @('833\u00f76=138, 5', '236\u00f78=29, 4'),\n  @('793\u00f72=396, 1', '508\u00f72=254, 0'),\n  @('150\u00f72=75, 0', '899\u00f74=224, 3'),\n  @('116\u00f77=16, 4', '482\u00f76=80, 2'),\n  @('916\u00f79=101, 7', '437\u00f77=62, 3'),\n  @('479\u00f76=79, 5', '291\u00f73=97, 0'),\n  @('651\u00f72=325, 1', '299\u00f74=74, 3'),\n  @('530\u00f79=58, 8', '667\u00f74=166, 3'),\n  @('130\u00f74=32, 2', '292\u00f76=48, 4'),\n  @('889\u00f78=111, 1', '203\u00f72=101, 1'),\n  @('331\u00f73=110, 1', '882\u00f78=110, 2'),\n  @('523\u00f75=104, 3', '731\u00f77=104, 3'),\n  @('519\u00f73=173, 0', '589\u00f77=84, 1'),\n  @('550\u00f74=137, 2', '201\u00f77=28, 5'),\n  @('290\u00f75=58, 0', '957\u00f76=159, 3')\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $result = $find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceAll)\n  if (-not $result) {\n    throw \"Find/Replace failed for: $oldText\"\n  }\n}\n"}
